# Update for new layout
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update pin numbers (C/D columns) to reflect the new physical layout ---
$ws.Range("C10").Value = 180
$ws.Range("D10").Value = 188
$ws.Range("C11").Value = 181
$ws.Range("D11").Value = 189
$ws.Range("C12").Value = 182
$ws.Range("D12").Value = 190
$ws.Range("C13").Value = 183
$ws.Range("D13").Value = 191
$ws.Range("C14").Value = 184
$ws.Range("D14").Value = 192
$ws.Range("C15").Value = 185
$ws.Range("D15").Value = 193
$ws.Range("C16").Value = 186
$ws.Range("D16").Value = 194
$ws.Range("C17").Value = 187
$ws.Range("D17").Value = 195
$ws.Range("C18").Value = 196
$ws.Range("D18").Value = 204
$ws.Range("C19").Value = 197
$ws.Range("D19").Value = 205
$ws.Range("C20").Value = 198
$ws.Range("D20").Value = 206
$ws.Range("C21").Value = 199
$ws.Range("D21").Value = 207
$ws.Range("C22").Value = 200
$ws.Range("D22").Value = 208
$ws.Range("C23").Value = 201
$ws.Range("D23").Value = 209

# --- 2. Remove the now-unused placeholder row (old row 24) ---
$ws.Rows.Item(24).Delete()
$ws.Range("B24:B25").ClearContents()

# --- 3. Add the new "Macro" helper column (M) with ROKUHANPOINT() builder ---
$ws.Range("M1").Value = "Macro"
$ws.Range("A10").Value = "ROKUHANPOINT("
$ws.Range("M2:M23").Formula = '=CONCATENATE($A$10,B2,$A$5,E2,$A$5,F2,$A$5,"",$A$5,"POINT",B2,$A$6)'
$ws.Range("M1:M23").EntireColumn.AutoFit()

# --- 4. Update the active selection ---
$ws.Range("H15").Select()
